{"js": "// Locate the paragraph that contains the \"poweroff\" command text \u2014 this is\n// the paragraph that gets split into four paragraphs by the commit:\n//   1) \"Rex change it 1\"\n//   2) (empty paragraph holding the _GoBack bookmark)\n//   3) the original \"1# poweroff ... the machine.\" text, now wrapped with\n//      gramStart/gramEnd proofErr markers around \"1# \"\n//   4) \"Rex change it 1\"\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"poweroff\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not locate the 'poweroff' paragraph\");\n}\n\n// The four replacement paragraphs, expressed as raw WordprocessingML so we\n// can place the <w:proofErr> grammar markers exactly where Word put them\n// and keep the lone bookmark on its own empty paragraph.\nconst replacementBodyXml =\n  \"<w:p><w:r><w:t>Rex change it 1</w:t></w:r></w:p>\" +\n  '<w:p><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>' +\n  \"<w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\">     </w:t></w:r>' +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  '<w:r><w:t xml:space=\"preserve\">1# </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  \"<w:r><w:t>poweroff</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\">        #</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  \"<w:r><w:t>poweroff</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> the machine.</w:t></w:r>' +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  \"</w:p>\" +\n  \"<w:p><w:r><w:t>Rex change it 1</w:t></w:r></w:p>\";\n\nconst flatOpcXml =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  replacementBodyXml +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\n// Insert the new paragraphs right after the target paragraph, then delete\n// the original paragraph \u2014 inserting via a collapsed \"replace\" range keeps\n// the paragraph(s) that precede the target intact (a direct whole-range\n// \"replace\" with multi-paragraph XML has been observed to eat the\n// preceding paragraph mark).\nconst insertionPoint = target.getRange(Word.RangeLocation.end);\ninsertionPoint.insertOoxml(flatOpcXml, Word.InsertLocation.after);\nawait context.sync();\n\ntarget.delete();\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that contains the \"poweroff\" command text -- this is\n# the paragraph that gets split into four paragraphs by the commit:\n#   1) \"Rex change it 1\"\n#   2) (empty paragraph holding the _GoBack bookmark)\n#   3) the original \"1# poweroff ... the machine.\" text, now wrapped with\n#      gramStart/gramEnd proofErr markers around \"1# \"\n#   4) \"Rex change it 1\"\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*poweroff*\") {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not locate the 'poweroff' paragraph\"\n}\n\n$wNs = \"xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'\"\n\n$replacementXml =\n    \"<w:p $wNs><w:r><w:t>Rex change it 1</w:t></w:r></w:p>\" +\n    \"<w:p $wNs><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>\" +\n    \"<w:p $wNs>\" +\n        \"<w:r><w:t xml:space='preserve'>     </w:t></w:r>\" +\n        \"<w:proofErr w:type='gramStart'/>\" +\n        \"<w:r><w:t xml:space='preserve'>1# </w:t></w:r>\" +\n        \"<w:proofErr w:type='spellStart'/>\" +\n        \"<w:r><w:t>poweroff</w:t></w:r>\" +\n        \"<w:proofErr w:type='spellEnd'/>\" +\n        \"<w:r><w:t xml:space='preserve'>        #</w:t></w:r>\" +\n        \"<w:proofErr w:type='spellStart'/>\" +\n        \"<w:r><w:t>poweroff</w:t></w:r>\" +\n        \"<w:proofErr w:type='spellEnd'/>\" +\n        \"<w:r><w:t xml:space='preserve'> the machine.</w:t></w:r>\" +\n        \"<w:proofErr w:type='gramEnd'/>\" +\n    \"</w:p>\" +\n    \"<w:p $wNs><w:r><w:t>Rex change it 1</w:t></w:r></w:p>\"\n\n# Range.InsertXML replaces the contents of the target range with the given\n# WordprocessingML -- this both removes the original paragraph and inserts\n# the four replacement paragraphs in one shot, keeping the preceding\n# paragraph (\"Using SSH:\") intact.\n$target.Range.InsertXML($replacementXml)\n"}
